# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# cryptocurrency rows whose quoted values moved in this run.
#
# Column D prices are stored as plain text (values such as "27.130.21" use
# dots as thousands separators and are not valid numbers). When a new price
# does parse as a normal decimal number (e.g. "307.08"), a leading apostrophe
# is used so Excel keeps storing it as literal text instead of silently
# converting it to a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.130.21"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.898.25"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'307.08"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "'0.5229"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.3800"
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("D9").Value = "'0.07283"
$ws.Range("E9").Value = "  +0.56%  "

$ws.Range("D10").Value = "'21.31"
$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("D11").Value = "'0.9057"
$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").Value = "'0.08196"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").Value = "1.891.85"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").Value = "'95.42"
$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").Value = "'5.350"
$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "'0.000008654"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").Value = "27.173.74"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").Value = "'5.118"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").Value = "2.117.89"
$ws.Range("E22").Value = "  -1.47%  "

$ws.Range("D23").Value = "'10.79"
$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("D24").Value = "'6.467"
$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").Value = "'2.334"
$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").Value = "'149.60"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("D27").Value = "'18.26"
$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("D28").Value = "'1.742"
$ws.Range("E28").Value = "  -0.77%  "

$ws.Range("D29").Value = "'115.36"
$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("D31").Value = "'4.864"
$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").Value = "'0.09234"
$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("D33").Value = "'0.05047"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").Value = "'0.7940"
$ws.Range("E34").Value = "  -2.71%  "

$ws.Range("D35").Value = "'1.223"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("D37").Value = "'3.384"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").Value = "'2.658"
$ws.Range("E38").Value = "  +3.84%  "

$ws.Range("D39").Value = "'0.5737"
$ws.Range("E39").Value = "  +0.98%  "

$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("D41").Value = "'1.080"
$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("D42").Value = "'9.026"
$ws.Range("E42").Value = "  +1.30%  "

$ws.Range("D43").Value = "'6.615"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").Value = "'116.27"
$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").Value = "'0.1517"
$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("D46").Value = "'0.4896"
$ws.Range("E46").Value = "  +1.67%  "

$ws.Range("E47").Value = "  +0.22%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").Value = "'1.640"
$ws.Range("E49").Value = "  +1.94%  "

$ws.Range("D50").Value = "'38.54"
$ws.Range("E50").Value = "  +3.04%  "

$ws.Range("D51").Value = "'64.06"
$ws.Range("E51").Value = "  +0.75%  "
